$d = $word.ActiveDocument

# Paragraph 1: title - change style from Title to Heading1 and update text
$d.Paragraphs.Item(1).Range.Style = "Heading1"
$d.Paragraphs.Item(1).Range.Text = "100 Beginner JavaScript Array Practice Problems"

$d.Paragraphs.Item(2).Range.Text = "1. Create an array of 5 numbers. Print only the numbers greater than 10."
$d.Paragraphs.Item(3).Range.Text = "2. Create an array of 4 animals. Print each animal in uppercase."
$d.Paragraphs.Item(4).Range.Text = "3. Create an array of 3 movies. Add another movie using .push() and print the updated array."
$d.Paragraphs.Item(5).Range.Text = "4. Create an array of 5 numbers. Print the square of each number."
$d.Paragraphs.Item(6).Range.Text = "5. Create an array of 4 songs. Print each song in the format: 'Song 1: [song name]'."
$d.Paragraphs.Item(7).Range.Text = "6. Create an array of fruits and remove the last one using .pop(). Print the updated array."
$d.Paragraphs.Item(8).Range.Text = "7. Create an array with 3 colors. Replace the second color with 'blue'. Print the array."
$d.Paragraphs.Item(9).Range.Text = "8. Create an array of 6 numbers. Count how many numbers are even."
$d.Paragraphs.Item(10).Range.Text = "9. Create an array of words. Print the length of each word."
$d.Paragraphs.Item(11).Range.Text = "10. Create an array of 5 numbers. Print the sum of all the numbers."
$d.Paragraphs.Item(12).Range.Text = "11. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(13).Range.Text = "12. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(14).Range.Text = "13. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(15).Range.Text = "14. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(16).Range.Text = "15. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(17).Range.Text = "16. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(18).Range.Text = "17. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(19).Range.Text = "18. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(20).Range.Text = "19. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(21).Range.Text = "20. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(22).Range.Text = "21. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(23).Range.Text = "22. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(24).Range.Text = "23. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(25).Range.Text = "24. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(26).Range.Text = "25. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(27).Range.Text = "26. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(28).Range.Text = "27. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(29).Range.Text = "28. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(30).Range.Text = "29. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(31).Range.Text = "30. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(32).Range.Text = "31. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(33).Range.Text = "32. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(34).Range.Text = "33. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(35).Range.Text = "34. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(36).Range.Text = "35. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(37).Range.Text = "36. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(38).Range.Text = "37. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(39).Range.Text = "38. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(40).Range.Text = "39. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(41).Range.Text = "40. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(42).Range.Text = "41. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(43).Range.Text = "42. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(44).Range.Text = "43. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(45).Range.Text = "44. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(46).Range.Text = "45. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(47).Range.Text = "46. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(48).Range.Text = "47. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(49).Range.Text = "48. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(50).Range.Text = "49. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(51).Range.Text = "50. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(52).Range.Text = "51. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(53).Range.Text = "52. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(54).Range.Text = "53. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(55).Range.Text = "54. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(56).Range.Text = "55. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(57).Range.Text = "56. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(58).Range.Text = "57. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(59).Range.Text = "58. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(60).Range.Text = "59. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(61).Range.Text = "60. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(62).Range.Text = "61. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(63).Range.Text = "62. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(64).Range.Text = "63. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(65).Range.Text = "64. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(66).Range.Text = "65. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(67).Range.Text = "66. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(68).Range.Text = "67. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(69).Range.Text = "68. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(70).Range.Text = "69. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(71).Range.Text = "70. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(72).Range.Text = "71. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(73).Range.Text = "72. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(74).Range.Text = "73. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(75).Range.Text = "74. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(76).Range.Text = "75. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(77).Range.Text = "76. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(78).Range.Text = "77. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(79).Range.Text = "78. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(80).Range.Text = "79. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(81).Range.Text = "80. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(82).Range.Text = "81. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(83).Range.Text = "82. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(84).Range.Text = "83. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(85).Range.Text = "84. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(86).Range.Text = "85. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(87).Range.Text = "86. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(88).Range.Text = "87. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(89).Range.Text = "88. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(90).Range.Text = "89. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(91).Range.Text = "90. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(92).Range.Text = "91. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(93).Range.Text = "92. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(94).Range.Text = "93. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(95).Range.Text = "94. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(96).Range.Text = "95. Create an array of 8 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(97).Range.Text = "96. Create an array of 3 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(98).Range.Text = "97. Create an array of 4 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(99).Range.Text = "98. Create an array of 5 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(100).Range.Text = "99. Create an array of 6 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
$d.Paragraphs.Item(101).Range.Text = "100. Create an array of 7 items. Perform a common array operation such as push, pop, indexing, or a loop to print or modify elements."
